# Updated the stats after release
# Append the latest data-collection run (2023-10-01 / serial 45201) as a new
# row to the bottom of the "Data" table, then move the "most recent value"
# highlight (Stars/Forks columns) from the old last row to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$lo = $ws.ListObjects.Item("Data")

# Remember the row the table currently ends on before growing it.
$prevLastRow = $lo.Range.Rows.Count + $lo.Range.Row - 1

# Grow the table by one row (ref + autoFilter + dimension all follow).
$lo.ListRows.Add() | Out-Null
$newLastRow = $prevLastRow + 1

# Seed the new row with the previous row's formulas/number formats so every
# column keeps its existing look (dates, integers, the two "delta" columns,
# and the calculated-column formulas) before the real values are written in.
$prevRowRange = $ws.Range("A" + $prevLastRow + ":AK" + $prevLastRow)
$newRowRange = $ws.Range("A" + $newLastRow + ":AK" + $newLastRow)
$prevRowRange.Copy()
$newRowRange.PasteSpecial(-4122) | Out-Null

# --- New row values (release snapshot taken after the latest commit) ---
$ws.Range("A" + $newLastRow).Value = 45201
$ws.Range("B" + $newLastRow).Value = 343
$ws.Range("C" + $newLastRow).Value = 129
$ws.Range("D" + $newLastRow).Value = 113
$ws.Range("E" + $newLastRow).Value = 280
$ws.Range("F" + $newLastRow).Value = 231
$ws.Range("G" + $newLastRow).Value = 5713
$ws.Range("H" + $newLastRow).Formula = "=Data[[#This Row],[LoC]]-G" + $prevLastRow
$ws.Range("I" + $newLastRow).Value = 6841
$ws.Range("J" + $newLastRow).Value = 1964
$ws.Range("K" + $newLastRow).Value = 553
$ws.Range("L" + $newLastRow).Value = 290
$ws.Range("M" + $newLastRow).Value = 143
$ws.Range("N" + $newLastRow).Value = 60
$ws.Range("O" + $newLastRow).Value = 16
$ws.Range("P" + $newLastRow).Formula = "=SUM(Data[[#This Row],[Shell]:[Bash]])"
$ws.Range("Q" + $newLastRow).Formula = "=Data[[#This Row],[Total]]-P" + $prevLastRow
$ws.Range("R" + $newLastRow).Value = 2123
$ws.Range("S" + $newLastRow).Value = 4505
$ws.Range("T" + $newLastRow).Value = 71027
$ws.Range("U" + $newLastRow).Value = 48715
$ws.Range("V" + $newLastRow).Value = 2
$ws.Range("W" + $newLastRow).Value = 1
$ws.Range("X" + $newLastRow).Value = 271
$ws.Range("Y" + $newLastRow).Formula = "=Data[[#This Row],[Open issues]]+Data[[#This Row],[Closed issues]]"
$ws.Range("Z" + $newLastRow).Value = 0
$ws.Range("AA" + $newLastRow).Value = 176
$ws.Range("AB" + $newLastRow).Formula = "=Data[[#This Row],[Open pull requests]]+Data[[#This Row],[Closed pull requests]]"
$ws.Range("AC" + $newLastRow).Value = 159
$ws.Range("AD" + $newLastRow).Value = 164
$ws.Range("AE" + $newLastRow).Value = 7
$ws.Range("AF" + $newLastRow).Value = 0
$ws.Range("AG" + $newLastRow).Value = 135
$ws.Range("AH" + $newLastRow).Value = 1053
$ws.Range("AI" + $newLastRow).Value = 7
$ws.Range("AJ" + $newLastRow).ClearContents()
$ws.Range("AK" + $newLastRow).Formula = "=SUM(Data[[#This Row],[Running]:[GH runs]])"

# The "Stars"/"Forks" highlight always sits on the newest row only - clear it
# from the row that used to be last now that the new row carries it instead.
$ws.Range("B" + $prevLastRow + ":C" + $prevLastRow).ClearFormats()

# Move the selection down onto the new row, matching where the author
# finished editing.
$ws.Activate()
$ws.Range("AJ" + $newLastRow).Select()
